$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 516, pushing existing rows 516:540 down to 517:541
$ws.Rows(516).Insert()

# Populate the new row 516 with values: descriptive columns copied from the
# (now shifted) row 517, price/volume columns are new data points.
$ws.Range("A516").Value = 10
$ws.Range("B516").Value = "Vega Modelo de Temuco"
$ws.Range("C516").Value = "La Araucanía"
$ws.Range("D516").Value = 45267
$ws.Range("E516").Value = 9
$ws.Range("F516").Value = "Fruta"
$ws.Range("G516").Value = 100102
$ws.Range("H516").Value = "Cítricos"
$ws.Range("I516").Value = 100102006
$ws.Range("J516").Value = "Pomelo"
$ws.Range("K516").Value = "Start Ruby"
$ws.Range("L516").Value = "Primera"
$ws.Range("M516").Value = 60
$ws.Range("N516").Value = 14000
$ws.Range("O516").Value = 14000
$ws.Range("P516").Value = 14000
$ws.Range("Q516").Value = "$/bandeja 15 kilos granel"
$ws.Range("R516").Value = "Región de O'Higgins"
$ws.Range("S516").Value = 933
$ws.Range("T516").Value = 15

# Ensure the date column keeps the date-number-format style used elsewhere in column D
$ws.Range("D516").NumberFormat = $ws.Range("D517").NumberFormat
